$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Heading swap.
#    "Integration Testing Skeleton:" (the heading before the first,
#    mostly-empty table) becomes "Unit" + an (empty) _GoBack bookmark
#    + " Testing Skeleton:" -- two runs instead of one.
# ------------------------------------------------------------------
$findAnchor = $d.Content.Find
$findAnchor.Execute("Add local points sum to Team", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorEnd = $findAnchor.Parent.End

$scopedRange1 = $d.Range($anchorEnd, $d.Content.End)
$find1 = $scopedRange1.Find
$find1.Execute("Integration Testing Skeleton:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1 = $find1.Parent
$rng1.Text = "Unit"

$insertPoint = $d.Range($rng1.End, $rng1.End)
$insertPoint.InsertAfter(" Testing Skeleton:")

$bmPoint = $d.Range($rng1.End, $rng1.End)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ------------------------------------------------------------------
# 2) Second heading: "Unit" (standalone run, right before the second
#    table) becomes "Integration".
# ------------------------------------------------------------------
$scopedRange2 = $d.Range($bmPoint.End, $d.Content.End)
$find2 = $scopedRange2.Find
$find2.Execute("Unit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2 = $find2.Parent
$rng2.Text = "Integration"

# ------------------------------------------------------------------
# 3) Second table (Login Button / Signup Button / ... / Exit Button):
#    - The old _GoBack bookmark that lived in the last cell ("Return
#      to game menu") must disappear -- it is re-created from scratch
#      by deleting and re-adding that row's content.
#    - A new row "Open the Website" / "Takes you to the homepage with
#      the play button" is appended after it.
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(2)
$lastRow = $tbl.Rows.Item($tbl.Rows.Count)
$lastRow.Delete()

$exitRow = $tbl.Rows.Add()
$exitRow.Cells.Item(1).Range.Text = "Exit Button"
$exitRow.Cells.Item(2).Range.Text = "Return to game menu"

$websiteRow = $tbl.Rows.Add()
$websiteRow.Cells.Item(1).Range.Text = "Open the Website"
$websiteRow.Cells.Item(2).Range.Text = "Takes you to the homepage with the play button"

Write-Output "done"
